# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values, replacing the previous Strike# figures in column G
$kValues = @{
    2  = 7
    3  = 3
    4  = 8
    5  = 6
    6  = 5
    7  = 4
    8  = 5
    9  = 1
    10 = 1
    11 = 3
    12 = 2
    13 = 3
    14 = 4
    15 = 4
    16 = 4
    17 = 3
    18 = 3
    19 = 6
    20 = 2
    21 = 3
    22 = 4
    23 = 6
    24 = 2
    25 = 4
    26 = 5
    27 = 8
    28 = 3
    29 = 2
    30 = 0
    31 = 4
    32 = 2
    33 = 4
    34 = 1
    35 = 1
    36 = 2
    37 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
